# capital_commitments.xlsx : trim bulk test rows down to a 2-row "fund 2"
# smoke-test fixture (mirrors switching the source bulk-files folder from
# "fund 1" to "fund 2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the 8 extra demo rows (rows 4-11); the blank template rows
#        below them (old 12-21) shift up to become rows 4-13. -------------
$ws.Rows("4:11").Delete()

# --- 2. Re-point the filter database defined name at the new, smaller
#        data range. -------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=CapitalCommitment!`$A`$1:`$S`$3"
    }
}

# --- 3. Row 2: first fund-2 commitment line -------------------------------
$ws.Range("A2").Value = "TSTF3"
$ws.Range("B2").Value = "Demo Fund 2"
$ws.Range("F2").Value = "TSTF21"
$ws.Range("Q2").Value = "TSTF3"
$ws.Range("L2").Clear()

# --- 4. Row 3: second fund-2 commitment line, now with both a folio- and
#        fund-currency committed amount (USD folio currency this time) ----
$ws.Range("A3").Value = "TSTF4"
$ws.Range("B3").Value = "Demo Fund 2"
$ws.Range("C3").Value = 12500000
$ws.Range("D3").Value = 1000000000
$ws.Range("E3").Value = "USD"
$ws.Range("F3").Value = "TSTF22"
$ws.Range("Q3").Value = "TSTF4"
$ws.Range("L3").Clear()

# P3 becomes a (blank) date-formatted cell, matching the K column styling.
$ws.Range("P3").NumberFormat = $ws.Range("K3").NumberFormat
$ws.Range("P3").Font.Name = "Arial"
$ws.Range("P3").Font.Size = 10

# --- 5. Leave the cursor parked on A2, like the saved file does. ---------
$ws.Range("A2").Select()
